$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("F2").Value = 1

$ws.Range("B3").Value = 1160
$ws.Range("C3").Value = 1260
$ws.Range("E3").Value = 1400
$ws.Range("F3").Value = 1520

# New values added in rows 4, 5, 8
$ws.Range("C4").Value = 1240
$ws.Range("F4").Value = 1500

$ws.Range("C5").Value = 1275
$ws.Range("E5").Value = 1440

$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 50

$ws.Range("E9").Value = 500
$ws.Range("F9").Value = 500

# Update selection (active cell) to M20
$ws.Range("M20").Select()
